$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2315.4443
$ws.Range("I15").Value = 2315.4443
$ws.Range("K15").Value = 6946.3329
$ws.Range("M15").Value = -6777.3329
$ws.Range("H19").Value = 1321.4117
$ws.Range("I19").Value = 1275.7273
$ws.Range("K19").Value = 1275.7273
$ws.Range("M19").Value = -1100.7273
$ws.Range("H33").Value = 221.47058
$ws.Range("I33").Value = 197.8125
$ws.Range("K33").Value = 197.8125
$ws.Range("M33").Value = 31.1875
$ws.Range("H92").Value = 1275.7273
$ws.Range("I92").Value = 1275.7273
$ws.Range("K92").Value = 1275.7273
$ws.Range("M92").Value = -27.72730000000001
$ws.Range("H132").Value = 1328.3636
$ws.Range("I132").Value = 1179.4445
$ws.Range("K132").Value = 3538.3335
$ws.Range("M132").Value = -1008.3335
$ws.Range("H133").Value = 75000
$ws.Range("J133").Value = 75000
$ws.Range("L133").Value = 75000
$ws.Range("N133").Value = -85120
$ws.Range("H137").Value = 5759.5
$ws.Range("I137").Value = 9222
$ws.Range("K137").Value = 27666
$ws.Range("M137").Value = -25116
$ws.Range("H138").Value = 3732.5898
$ws.Range("J138").Value = 4134.6787
$ws.Range("L138").Value = 12404.0361
$ws.Range("N138").Value = -22684.0361

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 869.7692
$ws.Range("J2").Value = 859.875
$ws.Range("L2").Value = 859.875
$ws.Range("N2").Value = -1085.875
$ws.Range("H97").Value = 611.5
$ws.Range("I97").Value = 611.5
$ws.Range("K97").Value = 611.5
$ws.Range("M97").Value = -115.5
$ws.Range("H110").Value = 2563.182
$ws.Range("I110").Value = 2423.3
$ws.Range("K110").Value = 2423.3
$ws.Range("M110").Value = -378.3000000000002
$ws.Range("H116").Value = 869.7692
$ws.Range("J116").Value = 859.875
$ws.Range("L116").Value = 859.875
$ws.Range("N116").Value = -5447.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 869.7692
$ws.Range("J3").Value = 859.875
$ws.Range("L3").Value = 859.875
$ws.Range("N3").Value = -1087.875
$ws.Range("H86").Value = 1416.5
$ws.Range("I86").Value = 1416.5
$ws.Range("K86").Value = 1416.5
$ws.Range("M86").Value = -293.5
$ws.Range("H89").Value = 1416.5
$ws.Range("I89").Value = 1416.5
$ws.Range("K89").Value = 7082.5
$ws.Range("M89").Value = -1466.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6658.9565
$ws.Range("I31").Value = 2717.1667
$ws.Range("J31").Value = 10959.091
$ws.Range("K31").Value = 2717.1667
$ws.Range("L31").Value = 10959.091
$ws.Range("M31").Value = -2422.1667
$ws.Range("N31").Value = -11549.091
$ws.Range("H34").Value = 6658.9565
$ws.Range("I34").Value = 2717.1667
$ws.Range("J34").Value = 10959.091
$ws.Range("K34").Value = 2717.1667
$ws.Range("L34").Value = 10959.091
$ws.Range("M34").Value = -2515.1667
$ws.Range("N34").Value = -11363.091
$ws.Range("H132").Value = 3824
$ws.Range("I132").Value = 3118.5715
$ws.Range("K132").Value = 9355.7145
$ws.Range("M132").Value = -6825.7145
$ws.Range("H133").Value = 124447
$ws.Range("J133").Value = 124447
$ws.Range("L133").Value = 124447
$ws.Range("N133").Value = -129507

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 364
$ws.Range("I14").Value = 364
$ws.Range("K14").Value = 1092
$ws.Range("M14").Value = -919
$ws.Range("H56").Value = 14999
$ws.Range("I56").Value = 14999
$ws.Range("K56").Value = 14999
$ws.Range("M56").Value = -14469
$ws.Range("H68").Value = 945.4666999999999
$ws.Range("J68").Value = 955.8570999999999
$ws.Range("L68").Value = 2867.5713
$ws.Range("N68").Value = -4489.5713
$ws.Range("H71").Value = 945.4666999999999
$ws.Range("J71").Value = 955.8570999999999
$ws.Range("L71").Value = 8602.713899999999
$ws.Range("N71").Value = -16714.7139
$ws.Range("H92").Value = 50
$ws.Range("J92").Value = 50
$ws.Range("L92").Value = 150
$ws.Range("N92").Value = -2646
$ws.Range("H121").Value = 881.3333
$ws.Range("I121").Value = 526
$ws.Range("K121").Value = 1578
$ws.Range("M121").Value = -268
$ws.Range("H132").Value = 4145.3
$ws.Range("I132").Value = 2247
$ws.Range("K132").Value = 20223
$ws.Range("M132").Value = -17693
$ws.Range("H139").Value = 773.2
$ws.Range("I139").Value = 861.5
$ws.Range("J139").Value = 420
$ws.Range("K139").Value = 2584.5
$ws.Range("L139").Value = 1260
$ws.Range("M139").Value = 2555.5
$ws.Range("N139").Value = -11540

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 15379
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 15379
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 15379
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -15681
$ws.Range("H80").Value = 3669.45
$ws.Range("I80").Value = 2499.0908
$ws.Range("J80").Value = 5099.8887
$ws.Range("K80").Value = 2499.0908
$ws.Range("L80").Value = 5099.8887
$ws.Range("M80").Value = -1501.0908
$ws.Range("N80").Value = -7095.8887
$ws.Range("H83").Value = 3669.45
$ws.Range("I83").Value = 2499.0908
$ws.Range("J83").Value = 5099.8887
$ws.Range("K83").Value = 12495.454
$ws.Range("L83").Value = 25499.4435
$ws.Range("M83").Value = -7503.454
$ws.Range("N83").Value = -35483.4435
$ws.Range("J113").Value = 1000
$ws.Range("L113").Value = 1000
$ws.Range("N113").Value = -5340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3983.3333
$ws.Range("J22").Value = 3983.3333
$ws.Range("L22").Value = 3983.3333
$ws.Range("N22").Value = -4573.3333
$ws.Range("H27").Value = 3983.3333
$ws.Range("J27").Value = 3983.3333
$ws.Range("L27").Value = 3983.3333
$ws.Range("N27").Value = -4197.3333
$ws.Range("H61").Value = 4815.8335
$ws.Range("I61").Value = 4724.5
$ws.Range("J61").Value = 4998.5
$ws.Range("K61").Value = 4724.5
$ws.Range("L61").Value = 4998.5
$ws.Range("M61").Value = -4522.5
$ws.Range("N61").Value = -5402.5
$ws.Range("H113").Value = 4815.8335
$ws.Range("I113").Value = 4724.5
$ws.Range("J113").Value = 4998.5
$ws.Range("K113").Value = 4724.5
$ws.Range("L113").Value = 4998.5
$ws.Range("M113").Value = -2554.5
$ws.Range("N113").Value = -9338.5
$ws.Range("H132").Value = 4537.923
$ws.Range("I132").Value = 4471.4287
$ws.Range("J132").Value = 4615.5
$ws.Range("K132").Value = 13414.2861
$ws.Range("L132").Value = 13846.5
$ws.Range("M132").Value = -10884.2861
$ws.Range("N132").Value = -18906.5
$ws.Range("H136").Value = 26965.945
$ws.Range("I136").Value = 8539.5
$ws.Range("K136").Value = 25618.5
$ws.Range("M136").Value = -23068.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 1000
$ws.Range("J5").Value = 1000
$ws.Range("L5").Value = 1000
$ws.Range("N5").Value = -1224
$ws.Range("H126").Value = 1560.6
$ws.Range("I126").Value = 1537.1428
$ws.Range("K126").Value = 4611.428400000001
$ws.Range("M126").Value = -2141.428400000001
$ws.Range("H136").Value = 7647.7144
$ws.Range("I136").Value = 7647.7144
$ws.Range("K136").Value = 22943.1432
$ws.Range("M136").Value = -20393.1432
